# DEAS_Equipment.xlsx — "added a second echelon" edit
#
# Movement Arcs (sheet6) gains a second echelon: the old echelon-1 -> sink
# arcs (rows 6-7) become echelon-1 -> echelon-2 arcs (mirroring the
# echelon-0 -> echelon-1 block in rows 2-5), and new echelon-2 -> sink arcs
# are appended (rows 8-11). Storage Room Arcs / Event Room Arcs each gain a
# matching echelon-2 self-loop row, and Utility Arcs' sink self-loop moves
# from echelon 2 to echelon 3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Movement Arcs
# ---------------------------------------------------------------------
$wsMove = $wb.Worksheets.Item("Movement Arcs")

function Set-Row($ws, $r, $vals) {
    $ws.Cells.Item($r, 1).Value  = $vals[0]
    $ws.Cells.Item($r, 2).Value  = $vals[1]
    $ws.Cells.Item($r, 3).Value  = $vals[2]
    $ws.Cells.Item($r, 4).Value  = $vals[3]
    $ws.Cells.Item($r, 5).Value  = $vals[4]
    $ws.Cells.Item($r, 6).Value  = $vals[5]
    $ws.Cells.Item($r, 7).Value  = $vals[6]
    $ws.Cells.Item($r, 8).Value  = $vals[7]
    $ws.Cells.Item($r, 9).Value  = $vals[8]
    $ws.Cells.Item($r, 10).Value = $vals[9]
}

# Rows 6 & 7 used to route echelon-1 nodes straight to the sink ("t" @
# echelon 2). They now route to the new echelon-2 nodes instead (same
# pattern as rows 2-5, shifted one echelon over).
Set-Row $wsMove 6 @("E1", 1, "b", "E1", 2, "a", "CHAIRS", 0, 72, 0)
Set-Row $wsMove 7 @("E1", 1, "b", "S1", 2, "a", "CHAIRS", 0, 72, 7)

# New rows for the rest of the echelon-1 -> echelon-2 block.
Set-Row $wsMove 8 @("S1", 1, "b", "E1", 2, "a", "CHAIRS", 0, 72, 7)
Set-Row $wsMove 9 @("S1", 1, "b", "S1", 2, "a", "CHAIRS", 0, 72, 0)

# New rows: echelon-2 nodes now feed the sink ("t" @ echelon 3).
Set-Row $wsMove 10 @("S1", 2, "b", "t", 3, "a", "CHAIRS", 0, 72, 0)
Set-Row $wsMove 11 @("E1", 2, "b", "t", 3, "a", "CHAIRS", 0, 72, 0)

$wsMove.Activate()
$wsMove.Range("I9").Select()

# ---------------------------------------------------------------------
# Storage Room Arcs
# ---------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Item("Storage Room Arcs")
Set-Row $wsStorage 3 @("S1", 2, "a", "S1", 2, "b", "CHAIRS", 0, 72, 0)

$wsStorage.Activate()
$wsStorage.Range("E3").Select()

# ---------------------------------------------------------------------
# Event Room Arcs
# ---------------------------------------------------------------------
$wsEvent = $wb.Worksheets.Item("Event Room Arcs")
Set-Row $wsEvent 3 @("E1", 2, "a", "E1", 2, "b", "CHAIRS", 0, 0, 0)

$wsEvent.Activate()
$wsEvent.Range("J3").Select()

# ---------------------------------------------------------------------
# Utility Arcs — sink self-loop moves from echelon 2 to echelon 3
# ---------------------------------------------------------------------
$wsUtility = $wb.Worksheets.Item("Utility Arcs")
$wsUtility.Cells.Item(4, 2).Value = 3
$wsUtility.Cells.Item(4, 5).Value = 3

$wsUtility.Activate()
$wsUtility.Range("H2").Select()

# ---------------------------------------------------------------------
# Extend the Movement Arcs filter database to the new data range
# ---------------------------------------------------------------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Movement Arcs!_FilterDatabase") {
        $n.RefersTo = "='Movement Arcs'!`$A`$1:`$J`$11"
    }
}
